# Natmi following Dr Hou advice
#
# The LR-pairs (Inhba-Tgfbr3) result table originally only modelled two
# clusters (FAPs, sCs). Per Dr Hou's advice a third cluster (ECs) is now
# included, turning the 2x2 sending/target cluster grid into a 3x3 grid
# (9 data rows instead of 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Inhba -> Tgfbr3 -> ECs
$ws.Cells.Item(2,1).Value  = "ECs"
$ws.Cells.Item(2,2).Value  = "Inhba"
$ws.Cells.Item(2,3).Value  = "Tgfbr3"
$ws.Cells.Item(2,4).Value  = "ECs"
$ws.Cells.Item(2,5).Value  = 2
$ws.Cells.Item(2,6).Value  = 0.6666666666666666
$ws.Cells.Item(2,7).Value  = 4.123204333333334
$ws.Cells.Item(2,8).Value  = 12.369613
$ws.Cells.Item(2,9).Value  = 0.2909967288544799
$ws.Cells.Item(2,10).Value = 0.2909967288544799
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 19.86261233333333
$ws.Cells.Item(2,14).Value = 59.587837
$ws.Cells.Item(2,15).Value = 0.1710751304955294
$ws.Cells.Item(2,16).Value = 0.1710751304955294
$ws.Cells.Item(2,17).Value = 81.89760924412012
$ws.Cells.Item(2,18).Value = 737.078483197081
$ws.Cells.Item(2,19).Value = 0.04978230336255235
$ws.Cells.Item(2,20).Value = 0.04978230336255235

# Row 3: ECs -> Inhba -> Tgfbr3 -> FAPs
$ws.Cells.Item(3,1).Value  = "ECs"
$ws.Cells.Item(3,2).Value  = "Inhba"
$ws.Cells.Item(3,3).Value  = "Tgfbr3"
$ws.Cells.Item(3,4).Value  = "FAPs"
$ws.Cells.Item(3,5).Value  = 2
$ws.Cells.Item(3,6).Value  = 0.6666666666666666
$ws.Cells.Item(3,7).Value  = 4.123204333333334
$ws.Cells.Item(3,8).Value  = 12.369613
$ws.Cells.Item(3,9).Value  = 0.2909967288544799
$ws.Cells.Item(3,10).Value = 0.2909967288544799
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 69.67747766666666
$ws.Cells.Item(3,14).Value = 209.032433
$ws.Cells.Item(3,15).Value = 0.6001266794307873
$ws.Cells.Item(3,16).Value = 0.6001266794307873
$ws.Cells.Item(3,17).Value = 287.2944778509366
$ws.Cells.Item(3,18).Value = 2585.650300658429
$ws.Cells.Item(3,19).Value = 0.1746349006126602
$ws.Cells.Item(3,20).Value = 0.1746349006126602

# Row 4: ECs -> Inhba -> Tgfbr3 -> sCs
$ws.Cells.Item(4,1).Value  = "ECs"
$ws.Cells.Item(4,2).Value  = "Inhba"
$ws.Cells.Item(4,3).Value  = "Tgfbr3"
$ws.Cells.Item(4,4).Value  = "sCs"
$ws.Cells.Item(4,5).Value  = 2
$ws.Cells.Item(4,6).Value  = 0.6666666666666666
$ws.Cells.Item(4,7).Value  = 4.123204333333334
$ws.Cells.Item(4,8).Value  = 12.369613
$ws.Cells.Item(4,9).Value  = 0.2909967288544799
$ws.Cells.Item(4,10).Value = 0.2909967288544799
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 26.564526
$ws.Cells.Item(4,14).Value = 79.693578
$ws.Cells.Item(4,15).Value = 0.2287981900736832
$ws.Cells.Item(4,16).Value = 0.2287981900736832
$ws.Cells.Item(4,17).Value = 109.530968716146
$ws.Cells.Item(4,18).Value = 985.7787184453141
$ws.Cells.Item(4,19).Value = 0.06657952487926735
$ws.Cells.Item(4,20).Value = 0.06657952487926735

# Row 5: FAPs -> Inhba -> Tgfbr3 -> ECs
$ws.Cells.Item(5,1).Value  = "FAPs"
$ws.Cells.Item(5,2).Value  = "Inhba"
$ws.Cells.Item(5,3).Value  = "Tgfbr3"
$ws.Cells.Item(5,4).Value  = "ECs"
$ws.Cells.Item(5,5).Value  = 3
$ws.Cells.Item(5,6).Value  = 1
$ws.Cells.Item(5,7).Value  = 8.433639666666666
$ws.Cells.Item(5,8).Value  = 25.300919
$ws.Cells.Item(5,9).Value  = 0.5952073574179045
$ws.Cells.Item(5,10).Value = 0.5952073574179045
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 19.86261233333333
$ws.Cells.Item(5,14).Value = 59.587837
$ws.Cells.Item(5,15).Value = 0.1710751304955294
$ws.Cells.Item(5,16).Value = 0.1710751304955294
$ws.Cells.Item(5,17).Value = 167.5141152580226
$ws.Cells.Item(5,18).Value = 1507.627037322203
$ws.Cells.Item(5,19).Value = 0.1018251763421673
$ws.Cells.Item(5,20).Value = 0.1018251763421673

# Row 6: FAPs -> Inhba -> Tgfbr3 -> FAPs
$ws.Cells.Item(6,1).Value  = "FAPs"
$ws.Cells.Item(6,2).Value  = "Inhba"
$ws.Cells.Item(6,3).Value  = "Tgfbr3"
$ws.Cells.Item(6,4).Value  = "FAPs"
$ws.Cells.Item(6,5).Value  = 3
$ws.Cells.Item(6,6).Value  = 1
$ws.Cells.Item(6,7).Value  = 8.433639666666666
$ws.Cells.Item(6,8).Value  = 25.300919
$ws.Cells.Item(6,9).Value  = 0.5952073574179045
$ws.Cells.Item(6,10).Value = 0.5952073574179045
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 69.67747766666666
$ws.Cells.Item(6,14).Value = 209.032433
$ws.Cells.Item(6,15).Value = 0.6001266794307873
$ws.Cells.Item(6,16).Value = 0.6001266794307873
$ws.Cells.Item(6,17).Value = 587.6347395228807
$ws.Cells.Item(6,18).Value = 5288.712655705926
$ws.Cells.Item(6,19).Value = 0.3571998149799809
$ws.Cells.Item(6,20).Value = 0.3571998149799809

# Row 7: FAPs -> Inhba -> Tgfbr3 -> sCs
$ws.Cells.Item(7,1).Value  = "FAPs"
$ws.Cells.Item(7,2).Value  = "Inhba"
$ws.Cells.Item(7,3).Value  = "Tgfbr3"
$ws.Cells.Item(7,4).Value  = "sCs"
$ws.Cells.Item(7,5).Value  = 3
$ws.Cells.Item(7,6).Value  = 1
$ws.Cells.Item(7,7).Value  = 8.433639666666666
$ws.Cells.Item(7,8).Value  = 25.300919
$ws.Cells.Item(7,9).Value  = 0.5952073574179045
$ws.Cells.Item(7,10).Value = 0.5952073574179045
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 26.564526
$ws.Cells.Item(7,14).Value = 79.693578
$ws.Cells.Item(7,15).Value = 0.2287981900736832
$ws.Cells.Item(7,16).Value = 0.2287981900736832
$ws.Cells.Item(7,17).Value = 224.035640199798
$ws.Cells.Item(7,18).Value = 2016.320761798182
$ws.Cells.Item(7,19).Value = 0.1361823660957564
$ws.Cells.Item(7,20).Value = 0.1361823660957564

# Row 8: sCs -> Inhba -> Tgfbr3 -> ECs
$ws.Cells.Item(8,1).Value  = "sCs"
$ws.Cells.Item(8,2).Value  = "Inhba"
$ws.Cells.Item(8,3).Value  = "Tgfbr3"
$ws.Cells.Item(8,4).Value  = "ECs"
$ws.Cells.Item(8,5).Value  = 3
$ws.Cells.Item(8,6).Value  = 1
$ws.Cells.Item(8,7).Value  = 1.612402333333333
$ws.Cells.Item(8,8).Value  = 4.837207
$ws.Cells.Item(8,9).Value  = 0.1137959137276156
$ws.Cells.Item(8,10).Value = 0.1137959137276156
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 19.86261233333333
$ws.Cells.Item(8,14).Value = 59.587837
$ws.Cells.Item(8,15).Value = 0.1710751304955294
$ws.Cells.Item(8,16).Value = 0.1710751304955294
$ws.Cells.Item(8,17).Value = 32.02652247236212
$ws.Cells.Item(8,18).Value = 288.238702251259
$ws.Cells.Item(8,19).Value = 0.01946765079080985
$ws.Cells.Item(8,20).Value = 0.01946765079080984

# Row 9: sCs -> Inhba -> Tgfbr3 -> FAPs
$ws.Cells.Item(9,1).Value  = "sCs"
$ws.Cells.Item(9,2).Value  = "Inhba"
$ws.Cells.Item(9,3).Value  = "Tgfbr3"
$ws.Cells.Item(9,4).Value  = "FAPs"
$ws.Cells.Item(9,5).Value  = 3
$ws.Cells.Item(9,6).Value  = 1
$ws.Cells.Item(9,7).Value  = 1.612402333333333
$ws.Cells.Item(9,8).Value  = 4.837207
$ws.Cells.Item(9,9).Value  = 0.1137959137276156
$ws.Cells.Item(9,10).Value = 0.1137959137276156
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 69.67747766666666
$ws.Cells.Item(9,14).Value = 209.032433
$ws.Cells.Item(9,15).Value = 0.6001266794307873
$ws.Cells.Item(9,16).Value = 0.6001266794307873
$ws.Cells.Item(9,17).Value = 112.3481275705146
$ws.Cells.Item(9,18).Value = 1011.133148134631
$ws.Cells.Item(9,19).Value = 0.0682919638381463
$ws.Cells.Item(9,20).Value = 0.06829196383814629

# Row 10: sCs -> Inhba -> Tgfbr3 -> sCs
$ws.Cells.Item(10,1).Value  = "sCs"
$ws.Cells.Item(10,2).Value  = "Inhba"
$ws.Cells.Item(10,3).Value  = "Tgfbr3"
$ws.Cells.Item(10,4).Value  = "sCs"
$ws.Cells.Item(10,5).Value  = 3
$ws.Cells.Item(10,6).Value  = 1
$ws.Cells.Item(10,7).Value  = 1.612402333333333
$ws.Cells.Item(10,8).Value  = 4.837207
$ws.Cells.Item(10,9).Value  = 0.1137959137276156
$ws.Cells.Item(10,10).Value = 0.1137959137276156
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 26.564526
$ws.Cells.Item(10,14).Value = 79.693578
$ws.Cells.Item(10,15).Value = 0.2287981900736832
$ws.Cells.Item(10,16).Value = 0.2287981900736832
$ws.Cells.Item(10,17).Value = 42.83270370629401
$ws.Cells.Item(10,18).Value = 385.4943333566461
$ws.Cells.Item(10,19).Value = 0.02603629909865945
$ws.Cells.Item(10,20).Value = 0.02603629909865945
